$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = 0.099974649717458

$ws.Range("B3").Value = 0.0870436897600398
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = 0.1870183394774978

$ws.Range("B4").Value = 0.1404651350714983
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 0.2404397847889563

$ws.Range("B5").Value = 0.04737453848398046
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = 0.1473491882014384

$ws.Range("B6").Value = 0.04164534289154693
$ws.Range("C6").Value = 0.00758821968299475
$ws.Range("D6").Value = 3.234991697728881
$ws.Range("E6").Value = 0.04027346932383607
$ws.Range("F6").Value = 0.02664642311026557
$ws.Range("G6").Value = 0.05664426267282812
$ws.Range("H6").Value = 0.1416199926090049

$ws.Range("B7").Value = 0.02766585445522624
$ws.Range("C7").Value = 0.00779535697115388
$ws.Range("D7").Value = 1.887231479864464
$ws.Range("E7").Value = 0.04059096086069674
$ws.Range("F7").Value = 0.01231001800819922
$ws.Range("G7").Value = 0.04302169090225323
$ws.Range("H7").Value = 0.1276405041726842

$ws.Range("B8").Value = 0.02470149967721591
$ws.Range("C8").Value = 0.004684855644768912
$ws.Range("D8").Value = 1.389282604641385
$ws.Range("E8").Value = 0.01315195166956533
$ws.Range("F8").Value = 0.01544761526769005
$ws.Range("G8").Value = 0.03395538408674198
$ws.Range("H8").Value = 0.1246761493946739

$ws.Range("B9").Value = 0.01888326388650609
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = 0.1188579136039641

$ws.Range("B10").Value = 0.02369878698745892
$ws.Range("C10").Value = 0.004273605237485141
$ws.Range("D10").Value = 1.042120459021991
$ws.Range("E10").Value = 0.01100731430426948
$ws.Range("F10").Value = 0.01519526873048314
$ws.Range("G10").Value = 0.03220230524443489
$ws.Range("H10").Value = 0.1236734367049169

$ws.Range("B11").Value = 0.02400530286068364
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = 0.1239799525781416

$ws.Range("B12").Value = 0.03883669499125091
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = 0.1388113447087089

$ws.Range("B13").Value = 0.04944992167375207
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("H13").Value = 0.1494245713912101

$ws.Range("B14").Value = 0.05444133214279879
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = ""
$ws.Range("H14").Value = 0.1544159818602568

$ws.Range("B15").Value = 0.06059297091547133
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
$ws.Range("G15").Value = ""
$ws.Range("H15").Value = 0.1605676206329293

$ws.Range("B16").Value = 0.06293349958891602
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = 0.162908149306374

$ws.Range("B17").Value = 0.06673435934246637
$ws.Range("C17").Value = 0.009175963386418757
$ws.Range("D17").Value = 11.51049643918773
$ws.Range("E17").Value = 0.05719195765868181
$ws.Range("F17").Value = 0.04864231096523789
$ws.Range("G17").Value = 0.08482640771969492
$ws.Range("H17").Value = 0.1667090090599244

$ws.Range("B18").Value = -0.099974649717458
$ws.Range("C18").Value = 0.01295305454255576
$ws.Range("D18").Value = -14.80351988215079
$ws.Range("E18").Value = 0.05011510350319052
$ws.Range("F18").Value = -0.1255028878472522
$ws.Range("G18").Value = -0.07444641158766359
$ws.Range("H18").Value = 0

$ws.Range("B19").Value = 0.06843837116250542
$ws.Range("C19").Value = 0.009476054681405473
$ws.Range("D19").Value = 12.02898695325092
$ws.Range("E19").Value = 0.05792681307000966
$ws.Range("F19").Value = 0.04973236581928257
$ws.Range("G19").Value = 0.08714437650572826
$ws.Range("H19").Value = 0.1684130208799634

$ws.Range("B20").Value = 0.0730291089370923
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = 0.1730037586545503

$ws.Range("B21").Value = 0.07102478244049937
$ws.Range("C21").Value = 0.009455517814283712
$ws.Range("D21").Value = 12.07657423541167
$ws.Range("E21").Value = 0.05720708924962065
$ws.Range("F21").Value = 0.05237497107718301
$ws.Range("G21").Value = 0.08967459380381548
$ws.Range("H21").Value = 0.1709994321579574

$ws.Range("B22").Value = 0.07417030871005445
$ws.Range("C22").Value = 0.009456611837933062
$ws.Range("D22").Value = 12.23124135819625
$ws.Range("E22").Value = 0.05186039735443365
$ws.Range("F22").Value = 0.05557909314717292
$ws.Range("G22").Value = 0.09276152427293605
$ws.Range("H22").Value = 0.1741449584275124

$ws.Range("B23").Value = 0.07723988582087125
$ws.Range("C23").Value = 0.00986551628199945
$ws.Range("D23").Value = 12.48328589402987
$ws.Range("E23").Value = 0.06059043561897473
$ws.Range("F23").Value = 0.05778157783675395
$ws.Range("G23").Value = 0.09669819380498867
$ws.Range("H23").Value = 0.1772145355383292

$ws.Range("B24").Value = 0.07718602813600076
$ws.Range("C24").Value = 0.009728480677344241
$ws.Range("D24").Value = 12.6064153140849
$ws.Range("E24").Value = 0.06213822187789601
$ws.Range("F24").Value = 0.05801727871975782
$ws.Range("G24").Value = 0.09635477755224384
$ws.Range("H24").Value = 0.1771606778534587

$ws.Range("B25").Value = 0.07734879572187313
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("G25").Value = ""
$ws.Range("H25").Value = 0.1773234454393311

$ws.Range("B26").Value = 0.07825374208207501
$ws.Range("C26").Value = 0.009748176530615457
$ws.Range("D26").Value = 12.50851141378558
$ws.Range("E26").Value = 0.05652768866642089
$ws.Range("F26").Value = 0.05901889424862341
$ws.Range("G26").Value = 0.09748858991552674
$ws.Range("H26").Value = 0.178228391799533

$ws.Range("B27").Value = 0.07918421373277997
$ws.Range("C27").Value = 0.009274613235989191
$ws.Range("D27").Value = 12.61071850318415
$ws.Range("E27").Value = 0.06203709638046995
$ws.Range("F27").Value = 0.06094615922916893
$ws.Range("G27").Value = 0.09742226823639089
$ws.Range("H27").Value = 0.179158863450238

$ws.Range("B28").Value = 0.08391664963040281
$ws.Range("C28").Value = 0.01020587934181395
$ws.Range("D28").Value = 12.5630685237041
$ws.Range("E28").Value = 0.09219881588401212
$ws.Range("F28").Value = 0.06382598427914278
$ws.Range("G28").Value = 0.1040073149816627
$ws.Range("H28").Value = 0.1838912993478608

$ws.Range("B29").Value = 0.02579560098016311
$ws.Range("C29").Value = 0.01163375223296176
$ws.Range("D29").Value = 3.930030570390517
$ws.Range("E29").Value = 0.008226197361249696
$ws.Range("F29").Value = 0.002764995551323798
$ws.Range("G29").Value = 0.04882620640900252
$ws.Range("H29").Value = 0.1257702506976211

